$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'243.34"
$ws.Range("D3").Formula = "'23.59"
$ws.Range("D4").Formula = "'5.286"
$ws.Range("D5").Formula = "'0.05780"
$ws.Range("D6").Formula = "'6.480"
$ws.Range("D7").Formula = "'3.331"
$ws.Range("D8").Formula = "'0.8084"
$ws.Range("D9").Formula = "'0.8785"
$ws.Range("D11").Formula = "'0.07276"
$ws.Range("D12").Formula = "'0.03090"
$ws.Range("D13").Formula = "'0.03058"
$ws.Range("D14").Formula = "'0.09315"
$ws.Range("D15").Formula = "'3.860"
$ws.Range("D16").Formula = "'0.001546"
$ws.Range("D17").Formula = "'0.04693"
$ws.Range("D18").Formula = "'0.0006048"
$ws.Range("D19").Formula = "'0.006045"
$ws.Range("D20").Formula = "'0.001294"
$ws.Range("D21").Formula = "'0.004602"
$ws.Range("D22").Formula = "'0.00008696"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("D23").Formula = "'3.580"
$ws.Range("D24").Formula = "'2.141"
$ws.Range("D41").Formula = "'0.006407"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Formula = "'0.003998"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Formula = "'0.1053"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Formula = "'0.007113"
$ws.Range("D45").Formula = "'0.00005468"
$ws.Range("D47").Formula = "'0.5498"
$ws.Range("D48").Formula = "'0.001858"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Formula = "'0.00002099"
$ws.Range("D50").Formula = "'0.0001999"

Write-Host "Applied 40 cell updates"
